# Edit: update match data for wales_cymru-premier_2023-2024.xlsx
# Applies: (1) reorder F:V content among existing rows 4-6, 8-10, 13-14, 15-17, 19-21
#          (2) append 4 new match rows (44-47)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rewrite F:V for rows that were reordered ---
# Row 4
$ws.Range("F4").Value = "Newtown"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Penybont"
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 2.92
$ws.Range("K4").Value = "11/08/2023 03:42"
$ws.Range("L4").Value = 4.03
$ws.Range("M4").Value = "12/08/2023 15:29"
$ws.Range("N4").Value = 3.27
$ws.Range("O4").Value = "11/08/2023 03:42"
$ws.Range("P4").Value = 3.34
$ws.Range("Q4").Value = "12/08/2023 15:26"
$ws.Range("R4").Value = 2.21
$ws.Range("S4").Value = "11/08/2023 03:42"
$ws.Range("T4").Value = 1.96
$ws.Range("U4").Value = "12/08/2023 15:26"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/wales/cymru-premier/newtown-penybont/0v7UYhSb/"

# Row 5
$ws.Range("F5").Value = "Bala"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "Barry"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1.29
$ws.Range("K5").Value = "11/08/2023 16:13"
$ws.Range("L5").Value = 1.38
$ws.Range("M5").Value = "12/08/2023 15:21"
$ws.Range("N5").Value = 5.49
$ws.Range("O5").Value = "11/08/2023 16:13"
$ws.Range("P5").Value = 4.94
$ws.Range("Q5").Value = "12/08/2023 15:21"
$ws.Range("R5").Value = 8.48
$ws.Range("S5").Value = "11/08/2023 16:13"
$ws.Range("T5").Value = 7.81
$ws.Range("U5").Value = "12/08/2023 15:21"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/wales/cymru-premier/bala-barry-town/Mm8QZYCh/"

# Row 6
$ws.Range("F6").Value = "Aberystwyth"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "Cardiff Metropolitan"
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3.08
$ws.Range("K6").Value = "11/08/2023 03:42"
$ws.Range("L6").Value = 3.41
$ws.Range("M6").Value = "12/08/2023 15:29"
$ws.Range("N6").Value = 3.29
$ws.Range("O6").Value = "11/08/2023 03:42"
$ws.Range("P6").Value = 3.29
$ws.Range("Q6").Value = "12/08/2023 15:29"
$ws.Range("R6").Value = 2.16
$ws.Range("S6").Value = "11/08/2023 03:42"
$ws.Range("T6").Value = 2.17
$ws.Range("U6").Value = "12/08/2023 15:21"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/wales/cymru-premier/aberystwyth-cardiff-metropolitan-university/hE3MzFcn/"

# Row 8
$ws.Range("F8").Value = "Connahs Q."
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = "Aberystwyth"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1.56
$ws.Range("K8").Value = "17/08/2023 09:13"
$ws.Range("L8").Value = 1.48
$ws.Range("M8").Value = "18/08/2023 20:36"
$ws.Range("N8").Value = 3.99
$ws.Range("O8").Value = "17/08/2023 09:13"
$ws.Range("P8").Value = 4.54
$ws.Range("Q8").Value = "18/08/2023 20:36"
$ws.Range("R8").Value = 4.72
$ws.Range("S8").Value = "17/08/2023 09:13"
$ws.Range("T8").Value = 6.14
$ws.Range("U8").Value = "18/08/2023 20:36"
$ws.Range("V8").Value = "https://www.betexplorer.com/football/wales/cymru-premier/connahs-q-aberystwyth/jTasWjCH/"

# Row 9
$ws.Range("F9").Value = "Caernarfon"
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = "Newtown"
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2.35
$ws.Range("K9").Value = "17/08/2023 09:13"
$ws.Range("L9").Value = 2.54
$ws.Range("M9").Value = "18/08/2023 19:47"
$ws.Range("N9").Value = 3.32
$ws.Range("O9").Value = "17/08/2023 09:13"
$ws.Range("P9").Value = 3.57
$ws.Range("Q9").Value = "18/08/2023 19:47"
$ws.Range("R9").Value = 2.67
$ws.Range("S9").Value = "17/08/2023 09:13"
$ws.Range("T9").Value = 2.61
$ws.Range("U9").Value = "18/08/2023 19:47"
$ws.Range("V9").Value = "https://www.betexplorer.com/football/wales/cymru-premier/caernarfon-newtown/d60wXWdB/"

# Row 10
$ws.Range("F10").Value = "Penybont"
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = "Haverfordwest"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1.72
$ws.Range("K10").Value = "17/08/2023 09:13"
$ws.Range("L10").Value = 1.46
$ws.Range("M10").Value = "18/08/2023 20:36"
$ws.Range("N10").Value = 3.56
$ws.Range("O10").Value = "17/08/2023 09:13"
$ws.Range("P10").Value = 4.22
$ws.Range("Q10").Value = "18/08/2023 20:36"
$ws.Range("R10").Value = 4.17
$ws.Range("S10").Value = "17/08/2023 09:13"
$ws.Range("T10").Value = 7.45
$ws.Range("U10").Value = "18/08/2023 20:36"
$ws.Range("V10").Value = "https://www.betexplorer.com/football/wales/cymru-premier/penybont-haverfordwest/UJboVARN/"

# Row 13
$ws.Range("F13").Value = "Pontypridd"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "Penybont"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3.75
$ws.Range("K13").Value = "24/08/2023 12:41"
$ws.Range("L13").Value = 4.07
$ws.Range("M13").Value = "25/08/2023 20:38"
$ws.Range("N13").Value = 3.44
$ws.Range("O13").Value = "24/08/2023 12:41"
$ws.Range("P13").Value = 3.21
$ws.Range("Q13").Value = "25/08/2023 20:37"
$ws.Range("R13").Value = 1.87
$ws.Range("S13").Value = "24/08/2023 12:41"
$ws.Range("T13").Value = 2
$ws.Range("U13").Value = "25/08/2023 20:38"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/wales/cymru-premier/pontypridd-united-penybont/M9q1nCdb/"

# Row 14
$ws.Range("F14").Value = "Bala"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "Connahs Q."
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2.72
$ws.Range("K14").Value = "24/08/2023 12:40"
$ws.Range("L14").Value = 2.68
$ws.Range("M14").Value = "25/08/2023 20:36"
$ws.Range("N14").Value = 3.07
$ws.Range("O14").Value = "24/08/2023 12:40"
$ws.Range("P14").Value = 3.41
$ws.Range("Q14").Value = "25/08/2023 20:36"
$ws.Range("R14").Value = 2.52
$ws.Range("S14").Value = "24/08/2023 12:40"
$ws.Range("T14").Value = 2.56
$ws.Range("U14").Value = "25/08/2023 20:36"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/wales/cymru-premier/bala-connahs-q/vofdmhsh/"

# Row 15
$ws.Range("F15").Value = "Haverfordwest"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "Caernarfon"
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 2.2
$ws.Range("K15").Value = "24/08/2023 23:42"
$ws.Range("L15").Value = 2.18
$ws.Range("M15").Value = "26/08/2023 15:23"
$ws.Range("N15").Value = 3.52
$ws.Range("O15").Value = "24/08/2023 23:42"
$ws.Range("P15").Value = 3.66
$ws.Range("Q15").Value = "26/08/2023 15:23"
$ws.Range("R15").Value = 2.84
$ws.Range("S15").Value = "24/08/2023 23:42"
$ws.Range("T15").Value = 3.07
$ws.Range("U15").Value = "26/08/2023 15:23"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/wales/cymru-premier/haverfordwest-caernarfon/UZn9pjRA/"

# Row 16
$ws.Range("F16").Value = "Cardiff Metropolitan"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Colwyn Bay"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1.69
$ws.Range("K16").Value = "24/08/2023 22:12"
$ws.Range("L16").Value = 1.75
$ws.Range("M16").Value = "26/08/2023 15:22"
$ws.Range("N16").Value = 3.63
$ws.Range("O16").Value = "24/08/2023 22:12"
$ws.Range("P16").Value = 3.62
$ws.Range("Q16").Value = "26/08/2023 15:22"
$ws.Range("R16").Value = 4.25
$ws.Range("S16").Value = "24/08/2023 22:12"
$ws.Range("T16").Value = 4.8
$ws.Range("U16").Value = "26/08/2023 15:22"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/wales/cymru-premier/cardiff-metropolitan-university-colwyn-bay/00r5oWB4/"

# Row 17
$ws.Range("F17").Value = "TNS"
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = "Barry"
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 1.03
$ws.Range("K17").Value = "24/08/2023 23:42"
$ws.Range("L17").Value = 1.03
$ws.Range("M17").Value = "26/08/2023 14:18"
$ws.Range("N17").Value = 20.77
$ws.Range("O17").Value = "24/08/2023 23:42"
$ws.Range("P17").Value = 27.47
$ws.Range("Q17").Value = "26/08/2023 15:05"
$ws.Range("R17").Value = 26.13
$ws.Range("S17").Value = "24/08/2023 23:42"
$ws.Range("T17").Value = 38.82
$ws.Range("U17").Value = "26/08/2023 15:05"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/wales/cymru-premier/tns-barry-town/jkyIrUeN/"

# Row 19
$ws.Range("F19").Value = "Connahs Q."
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = "Newtown"
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1.56
$ws.Range("K19").Value = "28/08/2023 09:12"
$ws.Range("L19").Value = 1.56
$ws.Range("M19").Value = "29/08/2023 20:37"
$ws.Range("N19").Value = 3.89
$ws.Range("O19").Value = "28/08/2023 09:12"
$ws.Range("P19").Value = 3.97
$ws.Range("Q19").Value = "29/08/2023 20:37"
$ws.Range("R19").Value = 4.84
$ws.Range("S19").Value = "28/08/2023 09:12"
$ws.Range("T19").Value = 6.18
$ws.Range("U19").Value = "29/08/2023 20:37"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/wales/cymru-premier/connahs-q-newtown/2sAypvhc/"

# Row 20
$ws.Range("F20").Value = "Colwyn Bay"
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = "TNS"
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 14.92
$ws.Range("K20").Value = "29/08/2023 11:12"
$ws.Range("L20").Value = 24.8
$ws.Range("M20").Value = "29/08/2023 20:39"
$ws.Range("N20").Value = 10.28
$ws.Range("O20").Value = "29/08/2023 11:12"
$ws.Range("P20").Value = 13.41
$ws.Range("Q20").Value = "29/08/2023 20:39"
$ws.Range("R20").Value = 1.11
$ws.Range("S20").Value = "29/08/2023 11:12"
$ws.Range("T20").Value = 1.08
$ws.Range("U20").Value = "29/08/2023 20:39"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/wales/cymru-premier/colwyn-bay-tns/OjBXpKxi/"

# Row 21
$ws.Range("F21").Value = "Penybont"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "Cardiff Metropolitan"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1.73
$ws.Range("K21").Value = "28/08/2023 09:12"
$ws.Range("L21").Value = 1.79
$ws.Range("M21").Value = "29/08/2023 20:36"
$ws.Range("N21").Value = 3.5
$ws.Range("O21").Value = "28/08/2023 09:12"
$ws.Range("P21").Value = 3.32
$ws.Range("Q21").Value = "29/08/2023 20:36"
$ws.Range("R21").Value = 4.21
$ws.Range("S21").Value = "28/08/2023 09:12"
$ws.Range("T21").Value = 5.09
$ws.Range("U21").Value = "29/08/2023 20:36"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/wales/cymru-premier/penybont-cardiff-metropolitan-university/d49uqb73/"

# --- Step 2: append new rows 44-47 ---
# Row 44
$ws.Range("A2").Copy($ws.Range("A44"))
$ws.Range("E2").Copy($ws.Range("E44"))
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "wales"
$ws.Range("C44").Value = "cymru-premier"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45192.64583333334
$ws.Range("F44").Value = "Barry"
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = "Bala"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 4.5
$ws.Range("K44").Value = "21/09/2023 12:43"
$ws.Range("L44").Value = 4.74
$ws.Range("M44").Value = "23/09/2023 15:27"
$ws.Range("N44").Value = 3.82
$ws.Range("O44").Value = "21/09/2023 12:43"
$ws.Range("P44").Value = 3.64
$ws.Range("Q44").Value = "23/09/2023 15:27"
$ws.Range("R44").Value = 1.62
$ws.Range("S44").Value = "21/09/2023 12:43"
$ws.Range("T44").Value = 1.75
$ws.Range("U44").Value = "23/09/2023 15:27"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/wales/cymru-premier/barry-town-bala/AZ9Yjn5s/"

# Row 45
$ws.Range("A2").Copy($ws.Range("A45"))
$ws.Range("E2").Copy($ws.Range("E45"))
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "wales"
$ws.Range("C45").Value = "cymru-premier"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("E45").Value = 45192.64583333334
$ws.Range("F45").Value = "Colwyn Bay"
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = "Newtown"
$ws.Range("I45").Value = 4
$ws.Range("J45").Value = 3.36
$ws.Range("K45").Value = "21/09/2023 12:43"
$ws.Range("L45").Value = 3.88
$ws.Range("M45").Value = "23/09/2023 15:28"
$ws.Range("N45").Value = 3.5
$ws.Range("O45").Value = "21/09/2023 12:43"
$ws.Range("P45").Value = 3.87
$ws.Range("Q45").Value = "23/09/2023 15:27"
$ws.Range("R45").Value = 1.97
$ws.Range("S45").Value = "21/09/2023 12:43"
$ws.Range("T45").Value = 1.85
$ws.Range("U45").Value = "23/09/2023 15:27"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/wales/cymru-premier/colwyn-bay-newtown/GnuZlQZg/"

# Row 46
$ws.Range("A2").Copy($ws.Range("A46"))
$ws.Range("E2").Copy($ws.Range("E46"))
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "wales"
$ws.Range("C46").Value = "cymru-premier"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45192.64583333334
$ws.Range("F46").Value = "Haverfordwest"
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = "Connahs Q."
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = 4.06
$ws.Range("K46").Value = "21/09/2023 12:43"
$ws.Range("L46").Value = 6.32
$ws.Range("M46").Value = "23/09/2023 15:25"
$ws.Range("N46").Value = 3.7
$ws.Range("O46").Value = "21/09/2023 12:43"
$ws.Range("P46").Value = 4.06
$ws.Range("Q46").Value = "23/09/2023 15:25"
$ws.Range("R46").Value = 1.74
$ws.Range("S46").Value = "21/09/2023 12:43"
$ws.Range("T46").Value = 1.53
$ws.Range("U46").Value = "23/09/2023 15:25"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/wales/cymru-premier/haverfordwest-connahs-q/A9vwlpk0/"

# Row 47
$ws.Range("A2").Copy($ws.Range("A47"))
$ws.Range("E2").Copy($ws.Range("E47"))
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "wales"
$ws.Range("C47").Value = "cymru-premier"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45192.76041666666
$ws.Range("F47").Value = "Caernarfon"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "Pontypridd"
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 2.1
$ws.Range("K47").Value = "21/09/2023 12:43"
$ws.Range("L47").Value = 2.28
$ws.Range("M47").Value = "23/09/2023 18:13"
$ws.Range("N47").Value = 3.26
$ws.Range("O47").Value = "21/09/2023 12:43"
$ws.Range("P47").Value = 3.52
$ws.Range("Q47").Value = "23/09/2023 18:07"
$ws.Range("R47").Value = 3.13
$ws.Range("S47").Value = "21/09/2023 12:43"
$ws.Range("T47").Value = 2.99
$ws.Range("U47").Value = "23/09/2023 18:13"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/wales/cymru-premier/caernarfon-pontypridd-united/xdtVk6Km/"
